$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): update "想去人数" (F column) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 505
$wsExpo.Range("F3").Value = 5986
$wsExpo.Range("F5").Value = 83

# Sheet "全部类型" (All Types): same rows duplicated, update matching cells
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 505
$wsAll.Range("F3").Value = 5986
$wsAll.Range("F6").Value = 83
